$d = $word.ActiveDocument

$end = $d.Content
$end.Collapse(0)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The team conducted scrum meetings every</w:t></w:r><w:r><w:t xml:space="preserve"> Friday, meeting for thirty minutes from 2-2:30. During these </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>thirty minute</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> meetings, instead of working on the project, we would bring up any concerns with the sprints, such as trouble with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jira</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, or any difficulties with the workload. </w:t></w:r><w:r><w:t>If there were any issues with project implementation, or a need for increased understanding on certain ideas necessary to their part of the sprint</w:t></w:r><w:r><w:t>, we discussed it at the scrum meetings</w:t></w:r><w:r><w:t xml:space="preserve">. It was in such meetings that we finalized due dates and planned the progress of the project. Additionally, we </w:t></w:r><w:r><w:t>used about</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>five minutes to provide</w:t></w:r><w:r><w:t xml:space="preserve"> each other</w:t></w:r><w:r><w:t xml:space="preserve"> with</w:t></w:r><w:r><w:t xml:space="preserve"> feedback on our work and communication with</w:t></w:r><w:r><w:t>in the group, and used this feedback to effectively improve our respective interpersonal skills</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@

[void]$end.InsertXML($xml)
